# Updates to Script 2 (not finished); copy-paste MetaboAnalyst script from Rstudio to test
$wb = $excel.ActiveWorkbook

# "Job to Run" sheet: change the job name in A2 and update the selection
$wsJob = $wb.Worksheets.Item("Job to Run")
$wsJob.Range("A2").Value = "Anid_HE_TJGIp11_pos_2018"

# "All" sheet: update the selection (no data change on this sheet)
$wsAll = $wb.Worksheets.Item("All")
$wsAll.Range("A2").Select()

# Leave "Job to Run" as the active/selected sheet with A7 selected,
# matching the final state recorded in the workbook.
$wsJob.Activate()
$wsJob.Range("A7").Select()
